# Composite string pattern changes
# Replace .NET-style "{0}", "{1}", ... composite format placeholders
# with printf/format style "%s" / "%d" / "%N:s" placeholders in the
# localized string table on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: Win32\Release\Format.exe.String.Unit1.SHello ---
$ws.Range("B11").Value = "Hello %s!"
$ws.Range("C11").Value = "Moi %s!"
$ws.Range("E11").Value = "Hallo %s!"
$ws.Range("F11").Value = "Hallo %s!"

# --- Row 12: Win32\Release\Format.exe.String.Unit1.SHello2 ---
$ws.Range("B12").Value = "Hello %s and %s!"
$ws.Range("C12").Value = "Moi %s ja %s!"
$ws.Range("D12").Value = "Moi %1:s ja %0:s!"
$ws.Range("E12").Value = "Hallo %s und %s!"
$ws.Range("F12").Value = "Hallo %s en %s!"
$ws.Range("G12").Value = "שלום %s %s!"

# --- Row 13: Win32\Release\Format.exe.String.Unit1.SCount ---
$ws.Range("B13").Value = "%s has %d cars"
$ws.Range("C13").Value = "%s:lla on %d autoa"
$ws.Range("E13").Value = "%s hat %d Autos"
$ws.Range("F13").Value = "%s heeft %d auto's"

# --- Row 14: Win32\Release\Format.exe.String.Unit1.SCount2 ---
$ws.Range("B14").Value = "%d cars will pick up %s and %s"
$ws.Range("C14").Value = "%d autoa hakee %s:in ja %s:in"
$ws.Range("E14").Value = "Autos nehmen %1:s und %2:s auf."
$ws.Range("F14").Value = "Auto's zullen %d ophalen %s en %s."
$ws.Range("G14").Value = "%d מכוניות יאספו %s ו-%s"

# --- Row 15: Win32\Release\Format.exe.String.Unit1.SDouble ---
$ws.Range("B15").Value = "%0:s swims and %0:s skis"
$ws.Range("C15").Value = "%0:s ui ja %0:s hiihtää"
$ws.Range("D15").Value = "%0:s ui ja %0:s hiihtää"
$ws.Range("F15").Value = "%0:s Zwemmen en %0:s ski's"
